# IPDO.xlsx edit script
# - Fills in rows 15-25 of "Tabela1" sheet by replicating the pattern of row 13/14
#   (full row) for most rows, and just column A for rows 17 and 22.
# - Fills in row 26 by replicating the pattern of row 5 (the first data row).
# - Removes the stray empty cell R14.
# - Adds a new second worksheet named "pag2", which becomes the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rows that become an exact copy of row 13 (" 31 Maio de 2016" pattern)
$fullRows = @(15, 16, 18, 19, 20, 21, 23, 24, 25)
foreach ($r in $fullRows) {
    $ws1.Range("A13:Q13").Copy() | Out-Null
    $ws1.Range("A$r`:Q$r").PasteSpecial() | Out-Null
    $ws1.Range("A$r").Style = "Normal"
    $ws1.Range("H$r").Style = "Normal"
    $ws1.Range("P$r").Style = "Normal"
}

# Rows that only get column A populated
$aOnlyRows = @(17, 22)
foreach ($r in $aOnlyRows) {
    $ws1.Range("A13").Copy() | Out-Null
    $ws1.Range("A$r").PasteSpecial() | Out-Null
    $ws1.Range("A$r").Style = "Normal"
}

# Row 26 becomes a copy of row 5 (the first data row, " 01 Maio de 2016")
$ws1.Range("A5:Q5").Copy() | Out-Null
$ws1.Range("A26:Q26").PasteSpecial() | Out-Null
$ws1.Range("A26").Style = "Normal"
$ws1.Range("H26").Style = "Normal"
$ws1.Range("P26").Style = "Normal"

$excel.CutCopyMode = 0

# Drop the stray empty cell left in R14
$ws1.Range("R14").ClearContents()

# Add the new "pag2" worksheet right after "Tabela1"; it becomes the active sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "pag2"
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 0.511811024 * 72
$ws2.PageSetup.RightMargin = 0.511811024 * 72
$ws2.PageSetup.TopMargin = 0.787401575 * 72
$ws2.PageSetup.BottomMargin = 0.787401575 * 72
$ws2.PageSetup.HeaderMargin = 0.31496062 * 72
$ws2.PageSetup.FooterMargin = 0.31496062 * 72
$ws2.Range("A1").Select() | Out-Null

Write-Host "Edit complete"
